$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy date number formatting from B6 onto B7 (keeps the existing built-in
# date numFmt instead of minting a duplicate custom one), then set values.
$ws.Range("B6").Copy($ws.Range("B7"))

$ws.Range("A7").Value = "Frost, Evan"
$ws.Range("B7").Value = 42774
$ws.Range("C7").Value = "???"
$ws.Range("D7").Value = "???"
$ws.Range("E7").Value = "Made basic water, sand, and birch tree, fall and spring/summer."

$ws.Range("A7:E7").RowHeight = 45

$ws.Range("E7").Select()
